$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3862.5
$ws.Range("I113").Value = 2668.3333
$ws.Range("J113").Value = 4579
$ws.Range("K113").Value = 2668.3333
$ws.Range("L113").Value = 4579
$ws.Range("M113").Value = 585.6667000000002
$ws.Range("N113").Value = -11087

$ws.Range("H116").Value = 5463.4546
$ws.Range("J116").Value = 6085.5713
$ws.Range("L116").Value = 6085.5713
$ws.Range("N116").Value = -12969.5713

$ws.Range("H132").Value = 4466251.5
$ws.Range("I132").Value = 4809726
$ws.Range("J132").Value = 1078
$ws.Range("K132").Value = 14429178
$ws.Range("L132").Value = 3234
$ws.Range("M132").Value = -14426648
$ws.Range("N132").Value = -8294

$ws.Range("H138").Value = 3023.0933
$ws.Range("I138").Value = 1759.8182
$ws.Range("J138").Value = 4816.129
$ws.Range("K138").Value = 5279.4546
$ws.Range("L138").Value = 14448.387
$ws.Range("M138").Value = -139.4546
$ws.Range("N138").Value = -24728.387

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 785.87177
$ws.Range("I74").Value = 772.29034
$ws.Range("J74").Value = 838.5
$ws.Range("K74").Value = 772.29034
$ws.Range("L74").Value = 838.5
$ws.Range("M74").Value = 101.70966
$ws.Range("N74").Value = -2586.5

$ws.Range("H77").Value = 785.87177
$ws.Range("I77").Value = 772.29034
$ws.Range("J77").Value = 838.5
$ws.Range("K77").Value = 3861.4517
$ws.Range("L77").Value = 4192.5
$ws.Range("M77").Value = 506.5482999999999
$ws.Range("N77").Value = -12928.5

$ws.Range("H132").Value = 2154.5386
$ws.Range("I132").Value = 1700.0333
$ws.Range("J132").Value = 3669.5557
$ws.Range("K132").Value = 5100.0999
$ws.Range("L132").Value = 11008.6671
$ws.Range("M132").Value = -2570.0999
$ws.Range("N132").Value = -16068.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1864.6818
$ws.Range("I20").Value = 1722.5714
$ws.Range("J20").Value = 2113.375
$ws.Range("K20").Value = 1722.5714
$ws.Range("L20").Value = 2113.375
$ws.Range("M20").Value = -1475.5714
$ws.Range("N20").Value = -2607.375

$ws.Range("H81").Value = 24370
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 24370
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 24370
$ws.Range("N81").Value = -26492
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 24370
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 24370
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 73110
$ws.Range("N84").Value = -83718
$ws.Range("M84").ClearContents()

$ws.Range("H86").Value = 2952.9678
$ws.Range("I86").Value = 2290.08
$ws.Range("J86").Value = 5715
$ws.Range("K86").Value = 2290.08
$ws.Range("L86").Value = 5715
$ws.Range("M86").Value = -1167.08
$ws.Range("N86").Value = -7961

$ws.Range("H89").Value = 2952.9678
$ws.Range("I89").Value = 2290.08
$ws.Range("J89").Value = 5715
$ws.Range("K89").Value = 11450.4
$ws.Range("L89").Value = 28575
$ws.Range("M89").Value = -5834.4
$ws.Range("N89").Value = -39807

$ws.Range("H99").Value = 2194160.5
$ws.Range("I99").Value = 3206214.5
$ws.Range("J99").Value = 1376.8334
$ws.Range("K99").Value = 3206214.5
$ws.Range("L99").Value = 1376.8334
$ws.Range("M99").Value = -3204716.5
$ws.Range("N99").Value = -4372.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1364636.4
$ws.Range("I6").Value = 1364636.4
$ws.Range("K6").Value = 1364636.4
$ws.Range("M6").Value = -1364523.4

$ws.Range("H31").Value = 2827.4883
$ws.Range("I31").Value = 2451.724
$ws.Range("J31").Value = 3605.8572
$ws.Range("K31").Value = 2451.724
$ws.Range("L31").Value = 3605.8572
$ws.Range("M31").Value = -2156.724
$ws.Range("N31").Value = -4195.8572

$ws.Range("H34").Value = 2827.4883
$ws.Range("I34").Value = 2451.724
$ws.Range("J34").Value = 3605.8572
$ws.Range("K34").Value = 2451.724
$ws.Range("L34").Value = 3605.8572
$ws.Range("M34").Value = -2249.724
$ws.Range("N34").Value = -4009.8572

$ws.Range("H74").Value = 28366.77
$ws.Range("J74").Value = 28366.77
$ws.Range("L74").Value = 28366.77
$ws.Range("N74").Value = -30114.77

$ws.Range("H77").Value = 28366.77
$ws.Range("J77").Value = 28366.77
$ws.Range("L77").Value = 85100.31
$ws.Range("N77").Value = -93836.31

$ws.Range("H94").Value = 167976.45
$ws.Range("J94").Value = 126506.5
$ws.Range("L94").Value = 126506.5
$ws.Range("N94").Value = -127408.5

$ws.Range("H132").Value = 1986.4231
$ws.Range("I132").Value = 1461.8
$ws.Range("J132").Value = 3735.1667
$ws.Range("K132").Value = 4385.4
$ws.Range("L132").Value = 11205.5001
$ws.Range("M132").Value = -1855.4
$ws.Range("N132").Value = -16265.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2722.3333
$ws.Range("J55").Value = 2723.4119
$ws.Range("L55").Value = 8170.2357
$ws.Range("N55").Value = -8524.235700000001

$ws.Range("H113").Value = 922.3077
$ws.Range("I113").Value = 1107.4546
$ws.Range("J113").Value = 682.7059
$ws.Range("K113").Value = 3322.3638
$ws.Range("L113").Value = 2048.1177
$ws.Range("M113").Value = -1152.3638
$ws.Range("N113").Value = -6388.117700000001

$ws.Range("H122").Value = 1228.2609
$ws.Range("I122").Value = 314
$ws.Range("J122").Value = 1628.25
$ws.Range("K122").Value = 2826
$ws.Range("L122").Value = 14654.25
$ws.Range("M122").Value = -376
$ws.Range("N122").Value = -19554.25

$ws.Range("H131").Value = 954.6585
$ws.Range("I131").Value = 772
$ws.Range("J131").Value = 980.0278
$ws.Range("K131").Value = 2316
$ws.Range("L131").Value = 2940.0834
$ws.Range("M131").Value = 2724
$ws.Range("N131").Value = -13020.0834

$ws.Range("H137").Value = 4905268.5
$ws.Range("J137").Value = 6176419
$ws.Range("L137").Value = 18529257
$ws.Range("N137").Value = -18539457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 23.666666
$ws.Range("I2").Value = 18.727272
$ws.Range("J2").Value = 37.25
$ws.Range("K2").Value = 18.727272
$ws.Range("L2").Value = 37.25
$ws.Range("M2").Value = 94.272728
$ws.Range("N2").Value = -263.25

$ws.Range("H70").Value = 3958
$ws.Range("I70").Value = 3944.3809
$ws.Range("J70").Value = 3989.7778
$ws.Range("K70").Value = 3944.3809
$ws.Range("L70").Value = 3989.7778
$ws.Range("M70").Value = -3674.3809
$ws.Range("N70").Value = -4529.7778

$ws.Range("H73").Value = 3958
$ws.Range("I73").Value = 3944.3809
$ws.Range("J73").Value = 3989.7778
$ws.Range("K73").Value = 3944.3809
$ws.Range("L73").Value = 3989.7778
$ws.Range("M73").Value = -3008.3809
$ws.Range("N73").Value = -5861.7778

$ws.Range("H132").Value = 1944.9062
$ws.Range("I132").Value = 1359.7894
$ws.Range("J132").Value = 2800.077
$ws.Range("K132").Value = 4079.3682
$ws.Range("L132").Value = 8400.231
$ws.Range("M132").Value = -1549.3682
$ws.Range("N132").Value = -13460.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3229.7083
$ws.Range("I122").Value = 3325.65
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 9976.950000000001
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -7526.950000000001
$ws.Range("N122").Value = -13150

$ws.Range("H132").Value = 4033.037
$ws.Range("I132").Value = 1835.5349
$ws.Range("J132").Value = 12623.272
$ws.Range("K132").Value = 5506.6047
$ws.Range("L132").Value = 37869.81600000001
$ws.Range("M132").Value = -2976.6047
$ws.Range("N132").Value = -42929.81600000001

$ws.Range("H136").Value = 2495.6985
$ws.Range("I136").Value = 1965.3112
$ws.Range("J136").Value = 3821.6667
$ws.Range("K136").Value = 5895.9336
$ws.Range("L136").Value = 11465.0001
$ws.Range("M136").Value = -3345.9336
$ws.Range("N136").Value = -16565.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 915.0164
$ws.Range("I132").Value = 809.875
$ws.Range("J132").Value = 1303.2307
$ws.Range("K132").Value = 2429.625
$ws.Range("L132").Value = 3909.6921
$ws.Range("M132").Value = 100.375
$ws.Range("N132").Value = -8969.6921
